$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for camera_id 5737 (row 191) was removed entirely; everything
# below shifts up by one (5738 -> row 191, 5739 -> row 192).
$ws.Rows(191).Delete()

# New hand-count values for camera_id 5664..5736 (rows 118..190).
$counts = @{
    118 = 2
    119 = 1
    120 = 1
    121 = 7
    122 = 0
    123 = 47
    124 = 2
    125 = 23
    126 = 2
    127 = 17
    128 = 1
    129 = 5
    130 = 0
    131 = 17
    132 = 0
    133 = 0
    134 = 8
    135 = 0
    136 = 0
    137 = 0
    138 = 9
    139 = 0
    140 = 5
    141 = 4
    142 = 3
    143 = 22
    144 = 0
    145 = 14
    146 = 0
    147 = 0
    148 = 27
    149 = 0
    150 = 25
    151 = 4
    152 = 1
    153 = 10
    154 = 5
    155 = 6
    156 = 1
    157 = 2
    158 = 18
    159 = 1
    160 = 0
    161 = 0
    162 = 0
    163 = 2
    164 = 7
    165 = 2
    166 = 1
    167 = 21
    168 = 0
    169 = 25
    170 = 5
    171 = 8
    172 = 0
    173 = 1
    174 = 2
    175 = 4
    176 = 10
    177 = 20
    178 = 0
    179 = 2
    180 = 0
    181 = 0
    182 = 1
    183 = 0
    184 = 1
    185 = 0
    186 = 0
    187 = 0
    188 = 2
    189 = 0
    190 = 0
}

foreach ($row in $counts.Keys) {
    $ws.Cells.Item($row, 2).Value = $counts[$row]
}

# Last two rows (camera_id 5738 and 5739, now at rows 191/192) get text
# annotations plus a stray measurement in C192.
$ws.Range("B191").Value = "NA"
$ws.Range("B192").Value = "cms"
$ws.Range("C192").Value = 118

# Scroll/selection state restored to match the saved view.
$ws.Range("C134").Select()
